{"js": "// The merge-field placeholder for the filer's/other party's address was\n// written with a typo (\".on_one_line\" instead of \".one_line\"), and in one\n// spot it was accidentally doubled up as \".address.address\" instead of\n// \".address.one_line\". Fix all occurrences so the address renders on one\n// line as intended (\"fixed issue of address not fitting on one line\").\n\nconst body = context.document.body;\n\n// 1) \"{{ users[0].address.address }}\" -> \"{{ users[0].address.one_line }}\"\nconst badDoubled = body.search(\"address.address\", { matchCase: true, matchWildcards: false });\nbadDoubled.load(\"items/text\");\n\n// 2) \"{{ *.address.on_one_line }}\" -> \"{{ *.address.one_line }}\"\nconst badTypo = body.search(\"address.on_one_line\", { matchCase: true, matchWildcards: false });\nbadTypo.load(\"items/text\");\n\nawait context.sync();\n\nfor (let i = 0; i < badDoubled.items.length; i++) {\n  badDoubled.items[i].insertText(\"address.one_line\", Word.InsertLocation.replace);\n}\n\nfor (let i = 0; i < badTypo.items.length; i++) {\n  badTypo.items[i].insertText(\"address.one_line\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The merge-field placeholder for the filer's/other party's address was\n# written with a typo (\".on_one_line\" instead of \".one_line\"), and in one\n# spot it was accidentally doubled up as \".address.address\" instead of\n# \".address.one_line\". Fix all occurrences so the address renders on one\n# line as intended (\"fixed issue of address not fitting on one line\").\n\n$d = $word.ActiveDocument\n\n# 1) \"{{ users[0].address.address }}\" -> \"{{ users[0].address.one_line }}\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"address.address\"\n$find1.Replacement.Text = \"address.one_line\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) \"{{ *.address.on_one_line }}\" -> \"{{ *.address.one_line }}\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"address.on_one_line\"\n$find2.Replacement.Text = \"address.one_line\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
